$d = $word.ActiveDocument

$replacements = @(
    @{old = "283÷9=31, 4"; new = "326÷4=81, 2"},
    @{old = "770÷8=96, 2"; new = "356÷4=89, 0"},
    @{old = "693÷3=231, 0"; new = "987÷4=246, 3"},
    @{old = "338÷6=56, 2"; new = "141÷2=70, 1"},
    @{old = "589÷9=65, 4"; new = "577÷8=72, 1"},
    @{old = "699÷6=116, 3"; new = "752÷5=150, 2"},
    @{old = "793÷4=198, 1"; new = "130÷6=21, 4"},
    @{old = "494÷9=54, 8"; new = "213÷5=42, 3"},
    @{old = "452÷6=75, 2"; new = "248÷5=49, 3"},
    @{old = "709÷9=78, 7"; new = "797÷4=199, 1"},
    @{old = "101÷8=12, 5"; new = "556÷5=111, 1"},
    @{old = "420÷3=140, 0"; new = "478÷8=59, 6"},
    @{old = "659÷7=94, 1"; new = "541÷5=108, 1"},
    @{old = "733÷5=146, 3"; new = "933÷8=116, 5"},
    @{old = "661÷5=132, 1"; new = "180÷2=90, 0"},
    @{old = "466÷2=233, 0"; new = "380÷6=63, 2"},
    @{old = "846÷8=105, 6"; new = "249÷2=124, 1"},
    @{old = "278÷3=92, 2"; new = "112÷2=56, 0"},
    @{old = "670÷2=335, 0"; new = "641÷7=91, 4"},
    @{old = "735÷7=105, 0"; new = "819÷4=204, 3"},
    @{old = "131÷5=26, 1"; new = "229÷9=25, 4"},
    @{old = "610÷5=122, 0"; new = "695÷2=347, 1"},
    @{old = "339÷3=113, 0"; new = "782÷2=391, 0"},
    @{old = "336÷2=168, 0"; new = "154÷7=22, 0"},
    @{old = "843÷2=421, 1"; new = "693÷4=173, 1"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
